$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Header for new column
$ws.Range("B1").Value = "Commission"

# Commission values for each subcategory row (2-14)
$ws.Range("B2").Value = 12
$ws.Range("B3").Value = 25
$ws.Range("B4").Value = 12
$ws.Range("B5").Value = 12
$ws.Range("B6").Value = 14
$ws.Range("B7").Value = 14
$ws.Range("B8").Value = 12
$ws.Range("B9").Value = 14
$ws.Range("B10").Value = 41
$ws.Range("B11").Value = 45
$ws.Range("B12").Value = 41
$ws.Range("B13").Value = 14
$ws.Range("B14").Value = 14

# Auto fit the new column width, matching the bestFit column behavior
$ws.Columns.Item(2).ColumnWidth = 11

# Update the selected cell, matching the diff's new selection
$ws.Range("C14").Select() | Out-Null
